$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the existing E-MTAB-7309-FILTERED database entry to GSE40279
$ws.Range("A2").Value = "GSE40279"

# Add two new database rows (EPIC, GSE55763) with the same begins/ends values
$ws.Range("A4").Value = "EPIC"
$ws.Range("B4").Value = 5
$ws.Range("C4").Value = 105

$ws.Range("A5").Value = "GSE55763"
$ws.Range("B5").Value = 5
$ws.Range("C5").Value = 105

# Widen column A to fit the longer entries
$ws.Columns.Item(1).ColumnWidth = 28.42578125

# Update the active selection
$null = $ws.Range("C11").Select()
